# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Xliff generate / handoff / handback timestamps are refreshed
# - Status/date columns are widened to fit the new "Ready for handoff" text

$wb = $excel.ActiveWorkbook

# Sheets: "Overview" (summary, cols A:G), "zh-cn" and "de-de" (detail, cols A:P)
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date (Overview!G2) and Latest Handback DateTime (de-de!H2):
# "2016-08-23 00:36:23" -> "2016-08-23 00:37:00"
$overview.Range("G2").Value = "2016-08-23 00:37:00"
$dede.Range("H2").Value = "2016-08-23 00:37:00"

# zh-cn Latest Handoff Datetime: "2016-08-23 00:36:18" -> "2016-08-23 00:36:55"
$zhcn.Range("H2").Value = "2016-08-23 00:36:55"

# Widen the Status/date columns to fit "Ready for handoff"
$overview.Range("E:F").ColumnWidth = 16.3
$zhcn.Range("C:C").ColumnWidth = 16.3
$dede.Range("C:C").ColumnWidth = 16.3
